$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (ano = 2025) metrics with refreshed values
$ws.Range("C8").Value = 1299
$ws.Range("D8").Value = 207
$ws.Range("E8").Value = 1092
$ws.Range("F8").Value = 8.490566037735849
$ws.Range("G8").Value = 84.06466512702079
$ws.Range("H8").Value = 15.93533487297922
